# Empleados.xlsx - fill in DESCUENTO (G) and SUELDO_QUIN (H) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1 header cell becomes the base quincenal salary value (replaces the
# "SUELDO_QUIN" label text) while keeping its existing style (s="1").
$ws.Cells.Item(1, 8).Value = 21000

# Row 26 (PEREZ MILLAN) also gets corrected ESTATUS/ASISTENCIAS/RETARDOS.
$ws.Cells.Item(26, 4).Value = "I"
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = 1

# Rows with 3+ RETARDOS get a 10% DESCUENTO and a discounted SUELDO_QUIN.
$discountRows = @(2, 7, 27, 30, 31, 37, 38, 39)
foreach ($r in $discountRows) {
    $cellG = $ws.Cells.Item($r, 7)
    $cellG.Value = 0.1
    $cellG.Style = "Normal"

    $cellH = $ws.Cells.Item($r, 8)
    $cellH.Value = 18900
    $cellH.Style = "Normal"
}

# All other data rows get the full SUELDO_QUIN (no descuento).
$fullRows = @(3, 4, 5, 6, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 28, 29, 32, 33, 34, 35, 36, 40, 41)
foreach ($r in $fullRows) {
    $cellH = $ws.Cells.Item($r, 8)
    $cellH.Value = 21000
    $cellH.Style = "Normal"
}
